$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.889.26"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "1.640.26"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.71%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.01"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5046"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2564"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06379"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.58"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07793"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.274"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.641.37"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5419"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").Value = "0.0₅7850"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.74"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").Value = "25.937.10"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.14"
$ws.Range("E19").Value = "  -2.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.378"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.947"
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.974"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.005"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.867"
$ws.Range("E24").Value = "  -3.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.10"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.841"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.71"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04932"
$ws.Range("E30").Value = "  -2.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.263"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.191"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.529"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.361"
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8928"
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.601"
$ws.Range("D37").Value = "1.139.44"
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5550"
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01559"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.005"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.670"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8190"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.35"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").Value = "0.0₈123"
$ws.Range("E44").Value = "  +8.93%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.778.74"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4511"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.005"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.11"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05056"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.007"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09488"
$ws.Range("E51").Value = "  +1.63%  "
